$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("__data")

for ($row = 6; $row -le 10; $row++) {
    $vfx = $ws.Range("L$row").Value2
    $sfx = $ws.Range("M$row").Value2
    $ws.Range("L$row").Value2 = "ui/assets/" + $vfx
    $ws.Range("M$row").Value2 = "ui/assets/" + $sfx
}
